# Mapeamento do balanço energético detalhado
# Duplicates the last data row (row 42) of sheet "Tabela1" into rows
# 43-48, extending the used range from A1:R42 to A1:R48.
#
# Copy/PasteSpecial (instead of assigning .Value) is used on purpose: it
# clones the source cells as-is (shared-string text cells), instead of
# letting the numeric-looking text ("27.136", "9.579", ...) get
# reinterpreted as numbers the way a plain .Value assignment would.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabela1")

for ($r = 43; $r -le 48; $r++) {
    $ws.Range("A42:Q42").Copy()
    $ws.Range("A" + $r + ":Q" + $r).PasteSpecial()
}

# Row 42's last column (R42) was an otherwise-empty text cell; row 48
# gets the same trailing placeholder cell.
$ws.Range("R42").Copy()
$ws.Range("R48").PasteSpecial()

$excel.CutCopyMode = $false
